# Updates symbol list values (price/volume columns + a handful of
# reshuffled coin rows) to match the 2023-01-06 10:52 UTC GitHub Actions
# data refresh. All Coin/Link/Price/Volume(1h) cells in this sheet are
# stored as literal text (inline strings, "General" format) rather than
# numbers/percentages, so we force text via NumberFormat "@" before
# assigning .Value, then restore the "Normal" style so the cell keeps
# its original (unstyled) look instead of picking up a Text number
# format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Range, [string]$Value) {
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "256.49"
Set-TextValue $ws.Range("E2") "-0.07%"
Set-TextValue $ws.Range("E3") "0.31%"
Set-TextValue $ws.Range("D4") "4.802"
Set-TextValue $ws.Range("E4") "1.83%"
Set-TextValue $ws.Range("D5") "0.05965"
Set-TextValue $ws.Range("E5") "0.58%"
Set-TextValue $ws.Range("D6") "6.638"
Set-TextValue $ws.Range("E6") "-0.36%"
Set-TextValue $ws.Range("D7") "0.8501"
Set-TextValue $ws.Range("E7") "-1.97%"
Set-TextValue $ws.Range("D8") "0.9258"
Set-TextValue $ws.Range("E8") "-1.60%"
Set-TextValue $ws.Range("D9") "0.1377"
Set-TextValue $ws.Range("E9") "-1.63%"
Set-TextValue $ws.Range("D10") "0.04231"
Set-TextValue $ws.Range("E10") "12.99%"
Set-TextValue $ws.Range("D11") "0.07011"
Set-TextValue $ws.Range("E11") "-2.03%"
Set-TextValue $ws.Range("D12") "0.03054"
Set-TextValue $ws.Range("E12") "-3.47%"
Set-TextValue $ws.Range("D13") "0.09090"
Set-TextValue $ws.Range("E13") "-1.70%"
Set-TextValue $ws.Range("D14") "0.001529"
Set-TextValue $ws.Range("E14") "-0.73%"
Set-TextValue $ws.Range("B15") "One"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D15") "0.0006067"
Set-TextValue $ws.Range("E15") "0.68%"
Set-TextValue $ws.Range("B16") "TigerCash"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D16") "0.006082"
Set-TextValue $ws.Range("E16") "0.66%"
Set-TextValue $ws.Range("B17") "LEO"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D17") "3.471"
Set-TextValue $ws.Range("E17") "-0.59%"
Set-TextValue $ws.Range("B18") "GateToken"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D18") "3.159"
Set-TextValue $ws.Range("E18") "-1.36%"
Set-TextValue $ws.Range("B19") "BTSEToken"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D19") "2.199"
Set-TextValue $ws.Range("E19") "-1.85%"
Set-TextValue $ws.Range("D20") "0.3041"
Set-TextValue $ws.Range("E20") "-2.73%"
Set-TextValue $ws.Range("E21") "0.31%"
Set-TextValue $ws.Range("D22") "3.904"
Set-TextValue $ws.Range("E22") "2.53%"
Set-TextValue $ws.Range("E23") "0.57%"
Set-TextValue $ws.Range("D24") "0.001219"
Set-TextValue $ws.Range("E24") "-0.14%"
Set-TextValue $ws.Range("D25") "0.003618"
Set-TextValue $ws.Range("E26") "0.04%"
Set-TextValue $ws.Range("E27") "1.92%"
Set-TextValue $ws.Range("D40") "0.03783"
Set-TextValue $ws.Range("E40") "-1.07%"
Set-TextValue $ws.Range("B41") "BKEXToken"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D41") "0.1099"
Set-TextValue $ws.Range("E41") "-0.33%"
Set-TextValue $ws.Range("B42") "KickToken"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D42") "0.003870"
Set-TextValue $ws.Range("E42") "-36.73%"
Set-TextValue $ws.Range("D43") "0.002449"
Set-TextValue $ws.Range("E43") "8.85%"
Set-TextValue $ws.Range("D44") "0.01414"
Set-TextValue $ws.Range("E44") "26.58%"
Set-TextValue $ws.Range("D45") "0.00005326"
Set-TextValue $ws.Range("E45") "-3.10%"
Set-TextValue $ws.Range("E46") "0.01%"
Set-TextValue $ws.Range("E47") "-50.29%"
Set-TextValue $ws.Range("D48") "0.3550"
Set-TextValue $ws.Range("E48") "14,646.37%"
Set-TextValue $ws.Range("E49") "0.01%"
Set-TextValue $ws.Range("E50") "0.01%"
